$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F holds "想去人数" (number interested)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 3398
$wsExhibition.Range("F4").Value = 133
$wsExhibition.Range("F5").Value = 6962
$wsExhibition.Range("F6").Value = 2410
$wsExhibition.Range("F7").Value = 39
$wsExhibition.Range("F8").Value = 107
$wsExhibition.Range("F14").Value = 567

# Sheet "全部类型" (All Types) - same column F meaning, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3398
$wsAll.Range("F5").Value = 133
$wsAll.Range("F6").Value = 6962
$wsAll.Range("F7").Value = 2410
$wsAll.Range("F8").Value = 39
$wsAll.Range("F9").Value = 107
$wsAll.Range("F15").Value = 567
